$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-24 01:55:42"

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
